$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Rename Sheet3 to Stats ---
$ws3.Name = "Stats"

# --- Sheet1 updates ---
$ws1.Range("D1").Value = "More"
$ws1.Range("B2").Value = "NoPage"
$ws1.Range("C2").Value = "Doména existuje, ale nejdou stáhnout žádné stránky, jinými slovy na webu neni žádný obsah."
$ws1.Range("B3").Value = "Parking"
$ws1.Range("C3").Value = "Když se jedná o doménu, která obsahuje webove stranky providera nebo se na strankách pracuje, nebo i forward na defaultní PARK stránku poskytovatele webu"
$ws1.Range("B4").Value = "Forwarding"
$ws1.Range("C4").Value = "Je tam forwarding mimo doménu, tj. forwarduje se na jinou doménu."
$ws1.Range("B5").Value = "Size"
$ws1.Range("C5").Value = "Když je velikost stránky příliš malá na kvalitní web, v tuto chvíli menší než 4000 bytů."
$ws1.Range("B6").Value = "GLinks"
$ws1.Range("C6").Value = "Obsahuje odkazy na oficialni web či adresu oficiálního vlastníka OZ"
$ws1.Range("B7").Value = "Title"
$ws1.Range("C7").Value = "Title obsahuje nazev OZ"
$ws1.Range("B8").Value = "GKeywords"
$ws1.Range("C8").Value = "keywords obsahuji  nazev OZ"
$ws1.Range("B9").Value = "SKeywords"
$ws1.Range("C9").Value = "podezrele keywords jako doména, seo, hosting, domény apod."
$ws1.Range("B10").Value = "Ads"
$ws1.Range("C10").Value = "Reklamy ruznych druhu..."
$ws1.Range("B11").Value = "SURL"
$ws1.Range("C11").Value = "Suspicios URL. Patern WWW?, obsahuje OZ ale k tomu spoustu dalsich znaku apod."
$ws1.Range("B12").Value = "Frames"
$ws1.Range("C12").Value = "Přiliš mnoho frames.. (více jak 2). Casto znaci reklamy... "
$ws1.Range("B13").Value = "Scontent"
$ws1.Range("C13").Value = "Uvodni stranka obsahuje OZ konkurence, sex, porno, kasino, slova, která nepatří do oboru Oz (třeba pro auta je to ubytování apod.)"
$ws1.Range("B14").Value = "Gowner"
$ws1.Range("C14").Value = "Dobrý vlastník domeny- owner domény obsahuje OZ"
$ws1.Range("B15").Value = "Bforward"
$ws1.Range("C15").Value = "Forward na stranky konkurence nebo na stranky, které jsou již označené jako BAD."
$ws1.Range("B16").Value = "SMS"
$ws1.Range("C16").Value = "Na strance je navedeni na poslani SMS apod."
$ws1.Range("B17").Value = "Blinks"
$ws1.Range("C17").Value = "Stranka obsahuje odkazy na stranky, oznacene jako BAD"
$ws1.Range("B18").Value = "NoOz"
$ws1.Range("C18").Value = " zatím není funkční - "
$ws1.Range("B19").Value = "Sowner"
$ws1.Range("C19").Value = "Vlastnik je proflaknutz podvodnik, tj. ma vice domen, ktere jsou marknute jako BAD autoklub-skoda.cz  clubskoda.cz skodafoto.cz "
$ws1.Range("B49").Value = "Generator"
$ws1.Range("C49").Value = "skoda-auto na auto-skoda wwwskoda peugeotlevne.cz peugeot-shop.cz tvujpeugeot.cz "
$ws1.Range("C50").Value = "myskoda peugeotcentrum.cz peugeot-auta, peugeot-auta peugeotweb.cz "
$ws1.Range("C51").Value = "e-skoda.cz eshopskoda.cz eskoda-shop.cz mercedesshop.com svetskoda nd-renault.cz renaultlevne.cz renaultnahradnidily renaultlaguna.cz renaultservis renaultweb"
$ws1.Range("D51").Value = "nahradnidilyhonda.cz  "
$ws1.Range("C52").Value = "info-skoda-auto.sk skoda-auto-web.cz skoda-club.cz servisskoda, servis-skoda skodadily skodahome skodalevne skodaleasing skodaslevy skodateam portalskoda skodaweb"
$ws1.Range("D52").Value = "hondaweb.cz "
$ws1.Range("B53").Value = "Toto maji byt stejni owneri"
$ws1.Range("C53").Value = "thalia-renault.cz  mercedescars.cz  mercedesforum.cz  mercedeslevne.cz  citroen-berlingo.cz  portal-citroen.cz  dilycitroen.cz  citroenstore.cz  citroenweb.cz  citroenbrno.cz"
$ws1.Range("C54").Value = "bmwlevne.cz bmw-olomouc.cz hyundaisantafe.cz   hyundaicz.cz hyundaiweb.cz  hyundaiix35.cz   hyundailevne.cz   portal-honda.cz   hondaclub.cz hondalevne.cz "
$ws1.Range("B55").Value = "SB:SUB000007124-ZONER Stanislav Skodak "
$ws1.Range("B56").Value = "SAVVY-1269883926 Stanislav Škodák      "
$ws1.Range("B57").Value = "oskodach.sk"

# --- Sheet2 content ---
$ws2.Range("A1").Value = "B1"
$ws2.Range("B1").Value = "Opravdovy zkudce - je u nej zamer na cizi IZ vydelat (parkuje, dava jako ze pekny obsah, ale ma reklamy apod)"
$ws2.Range("A2").Value = "B2"
$ws2.Range("B2").Value = "Maji svuj business, ale rozsirili nabidku o konkurenci, jedna se spise o omyl (lehke odstranit)"
$ws2.Range("A3").Value = "B3"
$ws2.Range("B3").Value = "Maji busness ale kompetne provazene (zameruji se na ruzne I konkurencni znacky a provazuji dokupy), tam je zamer prolinkovat, tezke odstranit"
$ws2.Range("B4").Value = "Yneuyiti pro vlastni propagaci (tedy ne zisk)"
$ws2.Range("B5").Value = "Skutencz vlastnik ma reklamu na konkurenci???"
$ws2.Range("B6").Value = "Privydelat reklamou roydil oproti cilenemu parazitovani"
$ws2.Range("B7").Value = "Co treba klubz< maji reklamu a tak je to v cajku<"

# --- Sheet3 content ---
$ws3.Range("A1").Value = "All domains"
$ws3.Range("B1").Value = 1873392
$ws3.Range("A2").Value = "CZ domains"
$ws3.Range("B2").Value = 568272
$ws3.Range("A3").Value = "CZ domains existed"
$ws3.Range("B3").Value = 467953
$ws3.Range("A5").Value = "All OZ"
$ws3.Range("B5").Value = 556969
$ws3.Range("A6").Value = "UPV"
$ws3.Range("B6").Value = 52479
$ws3.Range("A7").Value = "No duplicated UPV"
$ws3.Range("B7").Value = 48267
$ws3.Range("A8").Value = "OHIM"
$ws3.Range("B8").Value = 504490

# --- Sheet3: formulas for percentage columns ---
$ws3.Range("C2").Formula = "=B2/B1"
$ws3.Range("C3").Formula = "=B3/B2"
$ws3.Range("C6").Formula = "=B6/B5"
$ws3.Range("C7").Formula = "=B7/B5"
$ws3.Range("C8").Formula = "=B8/B5"

# --- Number formats ---
$ws3.Range("B1:B3").NumberFormat = "#\ ###\ ###"
$ws3.Range("B5:B8").NumberFormat = "#\ ###\ ###"
$ws3.Range("C2:C3").NumberFormat = "0%"
$ws3.Range("C6:C8").NumberFormat = "0%"

# --- Sheet3 column widths ---
$ws3.Columns.Item(1).ColumnWidth = 23.14
$ws3.Columns.Item(2).ColumnWidth = 38

# --- Sheet3 page setup (paper size 9 = A4, portrait) ---
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# --- Selections on each sheet ---
$ws1.Range("B15").Select()
$ws2.Range("A7").Select()

# --- Activate Stats (Sheet3) as active tab, with A7 selected ---
$ws3.Activate()
$ws3.Range("A7").Select()
